# BOM.xlsx minor update
#  - Row4 (Stainless steel Flat Washer): quantity note 28 -> 32 pcs
#  - Row5 (20 Series T Nuts): add quantity note "x8"
#  - Row6 (M3 M5 Hex Socket Bolts / Nuts): clarify size range in item name,
#    add quantity note "x20"
#  - Column A widened to fit the longer item text
#  - Selection moved to A7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Stainless steel Flat Washer - quantity x28 -> x32
$ws.Range("B4").Value = "x32"

# Row 5: 20 Series T Nuts M5 T Slot Nuts - add quantity x8
$ws.Range("B5").Value = "x8"

# Row 6: M3 M5 Hex Socket Bolts / Nuts - clarify size range, add quantity x20
$ws.Range("A6").Value = "M3 M5 Hex Socket Bolts / Nuts (5mm – 30mm)"
$ws.Range("B6").Value = "x20"

# Widen column A to fit the new, longer text
$ws.Columns("A").ColumnWidth = 63.67

# Move active selection to A7
$ws.Range("A7").Select()
